{"js": "// Edit 1: \"...how real application work on the real world \\n-there is some...\"\n//      -> \"...how real application work in the real world .\\n-there is some...\"\n// i.e. \"on\" -> \"in\" and the trailing space before the line break becomes \" .\"\n{\n  const results = context.document.body.search(\"work on the real world \", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"work in the real world .\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// Edit 2: \"2- familiarity with the delivery application  is medium risk:\"\n//      -> \"2- familiarity with the delivery System is medium risk:\"\n// The original text has an extra (invisible, RTL-tagged) space character between\n// \"application\" and \"is\" -- replacing \"application\" with \"System\" alone would leave\n// a double space, so we also remove that extra space to land on a single space,\n// matching the target text exactly.\n{\n  const anchor = context.document.body.search(\n    \"2- familiarity with the delivery application  is medium risk:\",\n    { matchCase: true, ignoreSpace: false }\n  );\n  anchor.load(\"items\");\n  await context.sync();\n  if (anchor.items.length > 0) {\n    const sentence = anchor.items[0];\n    // Split the sentence on single spaces (without trimming) so the lone extra\n    // space between \"application\" and \"is\" shows up as its own, isolated piece.\n    const pieces = sentence.split([\" \"], false, false);\n    pieces.load(\"items,text\");\n    await context.sync();\n\n    // Find \"application \" (word + its normal trailing space) and the lone\n    // single-space piece that immediately follows it.\n    let wordIdx = -1;\n    for (let i = 0; i < pieces.items.length - 1; i++) {\n      if (pieces.items[i].text === \"application \" && pieces.items[i + 1].text === \" \") {\n        wordIdx = i;\n        break;\n      }\n    }\n    if (wordIdx !== -1) {\n      // Delete the extra lone-space piece first (it comes after the word piece,\n      // so removing it does not invalidate the word piece's range).\n      pieces.items[wordIdx + 1].delete();\n      await context.sync();\n      pieces.items[wordIdx].insertText(\"System \", Word.InsertLocation.replace);\n      await context.sync();\n    }\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Edit 1 ---------------------------------------------------------------\n# \"...how real application work on the real world \\r-there is some...\"\n#  -> \"...how real application work in the real world .\\r-there is some...\"\n# \"on\" -> \"in\", and the trailing space before the line break becomes \" .\"\n$rng1 = $d.Content\n$found1 = $rng1.Find.Execute(\"work on the real world \")\nif ($found1) {\n    $rng1.Text = \"work in the real world .\"\n}\n\n# --- Edit 2 -----------------------------------------------------------------\n# \"2- familiarity with the delivery application  is medium risk:\"\n#  -> \"2- familiarity with the delivery System is medium risk:\"\n# Locate the unique occurrence of \"application\" right after \"the delivery \" in\n# that sentence (not the other \"application\" occurrences elsewhere in the\n# document), replace just that word, then remove the stray extra space\n# character that originally sat between \"application\" and \"is\".\n$rng2 = $d.Content\n$found2 = $rng2.Find.Execute(\"familiarity with the delivery application\")\nif ($found2) {\n    [void]$rng2.Collapse(0)\n    [void]$rng2.MoveStart(1, -11)\n    $rng2.Text = \"System\"\n\n    [void]$rng2.Collapse(0)\n    [void]$rng2.MoveEnd(1, 1)\n    if ($rng2.Text -eq \" \") {\n        $rng2.Delete()\n    }\n}\n"}
